$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.08692415871688995
$ws.Range("C2").Value = 0.4489496882726305
$ws.Range("D2").Value = 0.3262969795491077
$ws.Range("E2").Value = 0.5712241062394932
$ws.Range("F2").Value = 0.577260257759593

$ws.Range("B3").Value = 0.09922310024226133
$ws.Range("C3").Value = 0.6429000669399494
$ws.Range("D3").Value = 0.7520274682798707
$ws.Range("E3").Value = 0.8671951731184109
$ws.Range("F3").Value = 0.8808619649147968

$ws.Range("B4").Value = 0.5656050939845083
$ws.Range("C4").Value = 0.8647777313126734
$ws.Range("D4").Value = 4.153474585096863
$ws.Range("E4").Value = 2.038007503690029
$ws.Range("F4").Value = 2.001953655391513

$ws.Range("B5").Value = 0.2207140035010983
$ws.Range("C5").Value = 1.3294245198813
$ws.Range("D5").Value = 7.422634178464651
$ws.Range("E5").Value = 2.724451170137694
$ws.Range("F5").Value = 2.776526186899352

$ws.Range("B6").Value = 0.1296392257709591
$ws.Range("C6").Value = 1.189108630902367
$ws.Range("D6").Value = 7.156362825654766
$ws.Range("E6").Value = 2.675137907782469
$ws.Range("F6").Value = 2.73204717969972

$ws.Range("B7").Value = 0.2483310698430933
$ws.Range("C7").Value = 1.328958641558978
$ws.Range("D7").Value = 7.318429987019594
$ws.Range("E7").Value = 2.70525968938651
$ws.Range("F7").Value = 2.754380959454445

$ws.Range("B8").Value = 0.08227230788722427
$ws.Range("C8").Value = 1.318120586854214
$ws.Range("D8").Value = 7.426106611105771
$ws.Range("E8").Value = 2.72508836757742
$ws.Range("F8").Value = 2.785063824995722

$ws.Range("B9").Value = 0.1542305310830249
$ws.Range("C9").Value = 1.391611281455843
$ws.Range("D9").Value = 7.566127095767188
$ws.Range("E9").Value = 2.750659392903307
$ws.Range("F9").Value = 2.808055141407951

$ws.Range("B10").Value = 0.07646437456699881
$ws.Range("C10").Value = 1.325677262431856
$ws.Range("D10").Value = 7.411852158525697
$ws.Range("E10").Value = 2.722471700224944
$ws.Range("F10").Value = 2.782560325514475

$ws.Range("B11").Value = 0.09571397171271183
$ws.Range("C11").Value = 1.368091915780792
$ws.Range("D11").Value = 7.486065218487993
$ws.Range("E11").Value = 2.736067473306898
$ws.Range("F11").Value = 2.795847515386641
